# Auto-generated edit script: updates Leve profit calc columns (H-N)
# across 8 crafting-class sheets per scheduled market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5715.731
$ws.Range("I17").Value = 1350
$ws.Range("J17").Value = 6079.5415
$ws.Range("K17").Value = 4050
$ws.Range("L17").Value = 18238.6245
$ws.Range("M17").Value = -3882
$ws.Range("N17").Value = -18574.6245
$ws.Range("H32").Value = 10343.091
$ws.Range("I32").Value = 10732.833
$ws.Range("J32").Value = 9875.4
$ws.Range("K32").Value = 10732.833
$ws.Range("L32").Value = 9875.4
$ws.Range("M32").Value = -10406.833
$ws.Range("N32").Value = -10527.4
$ws.Range("H106").Value = 3103.074
$ws.Range("I106").Value = 1598.6666
$ws.Range("K106").Value = 1598.6666
$ws.Range("M106").Value = -967.6666
$ws.Range("H112").Value = 2257.6
$ws.Range("J112").Value = 2257.6
$ws.Range("L112").Value = 6772.799999999999
$ws.Range("N112").Value = -8988.799999999999
$ws.Range("H138").Value = 2952.1614
$ws.Range("I138").Value = 1907.5
$ws.Range("J138").Value = 3275.8591
$ws.Range("K138").Value = 5722.5
$ws.Range("L138").Value = 9827.577300000001
$ws.Range("M138").Value = -582.5
$ws.Range("N138").Value = -20107.5773

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3969.925
$ws.Range("I32").Value = 3558.923
$ws.Range("K32").Value = 3558.923
$ws.Range("M32").Value = -3271.923
$ws.Range("H42").Value = 22500
$ws.Range("I42").Value = 22500
$ws.Range("K42").Value = 22500
$ws.Range("M42").Value = -22014
$ws.Range("H74").Value = 66673270
$ws.Range("I74").Value = 111114680
$ws.Range("J74").Value = 11157
$ws.Range("K74").Value = 111114680
$ws.Range("L74").Value = 11157
$ws.Range("M74").Value = -111113806
$ws.Range("N74").Value = -12905
$ws.Range("H77").Value = 66673270
$ws.Range("I77").Value = 111114680
$ws.Range("J77").Value = 11157
$ws.Range("K77").Value = 555573400
$ws.Range("L77").Value = 55785
$ws.Range("M77").Value = -555569032
$ws.Range("N77").Value = -64521
$ws.Range("H110").Value = 5956.4
$ws.Range("I110").Value = 4324.85
$ws.Range("J110").Value = 12482.6
$ws.Range("K110").Value = 4324.85
$ws.Range("L110").Value = 12482.6
$ws.Range("M110").Value = -2279.85
$ws.Range("N110").Value = -16572.6
$ws.Range("H114").Value = 59398.6
$ws.Range("J114").Value = 59398.6
$ws.Range("L114").Value = 59398.6
$ws.Range("N114").Value = -68076.60000000001
$ws.Range("H122").Value = 5404.4375
$ws.Range("I122").Value = 5565.857
$ws.Range("J122").Value = 5278.8887
$ws.Range("K122").Value = 16697.571
$ws.Range("L122").Value = 15836.6661
$ws.Range("M122").Value = -14247.571
$ws.Range("N122").Value = -20736.6661
$ws.Range("H132").Value = 2201.5745
$ws.Range("I132").Value = 1496.4849
$ws.Range("K132").Value = 4489.4547
$ws.Range("M132").Value = -1959.4547
$ws.Range("H141").Value = 71963.336
$ws.Range("J141").Value = 71963.336
$ws.Range("L141").Value = 71963.336
$ws.Range("N141").Value = -82323.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3208.0833
$ws.Range("I134").Value = 1832
$ws.Range("K134").Value = 5496
$ws.Range("M134").Value = -2961
$ws.Range("H135").Value = 50844.184
$ws.Range("J135").Value = 50844.184
$ws.Range("L135").Value = 50844.184
$ws.Range("N135").Value = -60984.184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27209.143
$ws.Range("I31").Value = 1935.8857
$ws.Range("J31").Value = 153575.42
$ws.Range("K31").Value = 1935.8857
$ws.Range("L31").Value = 153575.42
$ws.Range("M31").Value = -1640.8857
$ws.Range("N31").Value = -154165.42
$ws.Range("H34").Value = 27209.143
$ws.Range("I34").Value = 1935.8857
$ws.Range("J34").Value = 153575.42
$ws.Range("K34").Value = 1935.8857
$ws.Range("L34").Value = 153575.42
$ws.Range("M34").Value = -1733.8857
$ws.Range("N34").Value = -153979.42
$ws.Range("H58").Value = 3685.577
$ws.Range("J58").Value = 5773
$ws.Range("L58").Value = 5773
$ws.Range("N58").Value = -6179
$ws.Range("H86").Value = 5153.222
$ws.Range("I86").Value = 4911.2856
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 4911.2856
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -3788.2856
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 5153.222
$ws.Range("I89").Value = 4911.2856
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 24556.428
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -18940.428
$ws.Range("N89").Value = -41232
$ws.Range("H99").Value = 2204.6
$ws.Range("I99").Value = 1760
$ws.Range("K99").Value = 1760
$ws.Range("M99").Value = -262
$ws.Range("H126").Value = 2204.6
$ws.Range("I126").Value = 1760
$ws.Range("K126").Value = 5280
$ws.Range("M126").Value = -2810
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H136").Value = 3685.577
$ws.Range("J136").Value = 5773
$ws.Range("L136").Value = 17319
$ws.Range("N136").Value = -22419
$ws.Range("H141").Value = 188168.58
$ws.Range("J141").Value = 243436.2
$ws.Range("L141").Value = 243436.2
$ws.Range("N141").Value = -253796.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3850408.5
$ws.Range("J131").Value = 2526970.5
$ws.Range("L131").Value = 7580911.5
$ws.Range("N131").Value = -7590991.5
$ws.Range("H136").Value = 4199.4
$ws.Range("I136").Value = 3999.25
$ws.Range("K136").Value = 11997.75
$ws.Range("M136").Value = -6897.75
$ws.Range("H137").Value = 61631.94
$ws.Range("J137").Value = 69696.2
$ws.Range("L137").Value = 209088.6
$ws.Range("N137").Value = -219288.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 14900
$ws.Range("J7").Value = 14900
$ws.Range("L7").Value = 14900
$ws.Range("N7").Value = -15124
$ws.Range("H8").Value = 14900
$ws.Range("J8").Value = 14900
$ws.Range("L8").Value = 14900
$ws.Range("N8").Value = -15178
$ws.Range("H80").Value = 561279
$ws.Range("I80").Value = 5000000
$ws.Range("J80").Value = 6438.875
$ws.Range("K80").Value = 5000000
$ws.Range("L80").Value = 6438.875
$ws.Range("M80").Value = -4999002
$ws.Range("N80").Value = -8434.875
$ws.Range("H83").Value = 561279
$ws.Range("I83").Value = 5000000
$ws.Range("J83").Value = 6438.875
$ws.Range("K83").Value = 25000000
$ws.Range("L83").Value = 32194.375
$ws.Range("M83").Value = -24995008
$ws.Range("N83").Value = -42178.375
$ws.Range("H122").Value = 11616.929
$ws.Range("I122").Value = 13382.0625
$ws.Range("J122").Value = 9263.416999999999
$ws.Range("K122").Value = 40146.1875
$ws.Range("L122").Value = 27790.251
$ws.Range("M122").Value = -37696.1875
$ws.Range("N122").Value = -32690.251
$ws.Range("H126").Value = 3334.4614
$ws.Range("I126").Value = 1705.9166
$ws.Range("J126").Value = 4730.357
$ws.Range("K126").Value = 5117.7498
$ws.Range("L126").Value = 14191.071
$ws.Range("M126").Value = -2647.7498
$ws.Range("N126").Value = -19131.071
$ws.Range("H138").Value = 65000
$ws.Range("J138").Value = 65000
$ws.Range("L138").Value = 65000
$ws.Range("N138").Value = -75280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 8710.333000000001
$ws.Range("I31").Value = 267.6
$ws.Range("K31").Value = 267.6
$ws.Range("M31").Value = -19.60000000000002
$ws.Range("H68").Value = 6550.8096
$ws.Range("I68").Value = 3599.5
$ws.Range("J68").Value = 7731.3335
$ws.Range("K68").Value = 3599.5
$ws.Range("L68").Value = 7731.3335
$ws.Range("M68").Value = -2850.5
$ws.Range("N68").Value = -9229.333500000001
$ws.Range("H71").Value = 6550.8096
$ws.Range("I71").Value = 3599.5
$ws.Range("J71").Value = 7731.3335
$ws.Range("K71").Value = 17997.5
$ws.Range("L71").Value = 38656.6675
$ws.Range("M71").Value = -14253.5
$ws.Range("N71").Value = -46144.6675
$ws.Range("H136").Value = 3035.9387
$ws.Range("I136").Value = 2180.4211
$ws.Range("J136").Value = 5991.364
$ws.Range("K136").Value = 6541.263300000001
$ws.Range("L136").Value = 17974.092
$ws.Range("M136").Value = -3991.263300000001
$ws.Range("N136").Value = -23074.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6278.5
$ws.Range("I126").Value = 6612.143
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 19836.429
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -17366.429
$ws.Range("N126").Value = -21440
$ws.Range("H136").Value = 4879.3335
$ws.Range("J136").Value = 6152.7856
$ws.Range("L136").Value = 18458.3568
$ws.Range("N136").Value = -23558.3568
